$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 150, shifting existing rows 150:234 down to 151:235.
$ws.Rows("150:150").Insert()

# Populate the freshly inserted row 150 with the new record. The other
# static columns (A,B,C,E,F,G,H,I,N,Q,R) match the surrounding rows, so
# copy them down from row 151 (the row that used to be row 150) to keep
# formatting/content consistent, then set the cells that actually carry
# new data for this record.
$ws.Range("A150").Value = 10
$ws.Range("B150").Value = 'Vega Modelo de Temuco'
$ws.Range("C150").Value = 'La Araucanía'
$ws.Range("D150").Value = 44960
$ws.Range("E150").Value = 9
$ws.Range("F150").Value = 100112012
$ws.Range("G150").Value = 'Espinaca'
$ws.Range("H150").Value = 'Sin especificar'
$ws.Range("I150").Value = 'Primera'
$ws.Range("J150").Value = 40
$ws.Range("K150").Value = 12000
$ws.Range("L150").Value = 12000
$ws.Range("M150").Value = 12000
$ws.Range("N150").Value = '$/docena de atados'
$ws.Range("O150").Value = 'Región de La Araucanía'
$ws.Range("P150").Value = 4000
$ws.Range("Q150").Value = 3
$ws.Range("R150").Value = 'Hortaliza'
